# Append/refresh scrape: 2025-10-07 18:27 JST
# Rebuilds the job-listing table (rows 2-22) on the "ランサーズ" sheet with the
# freshly scraped data, re-stamps the retrieval timestamp, and repoints the
# F-column hyperlinks so each row's link matches its URL text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$timestamp = "2025-10-07 18:27:06"

# Columns: B=title, C=category, D=price, E=deadline, F=url, G=score, H=skill summary (optional)
$rows = @(
    @{ B = "初回 【SES前提】業務委託エンジニア募集 自社AIサービス開発にも参画可能 業務システム開発の仕事の依頼"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408156"; G = 410; H = "🔥AI,Ai ◆開発,システム開発" },
    @{ B = "自社開発のロジシステムをサポート及び開発できる方募集【PHP, Python, VBA etc】"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408668"; G = 305; H = "🔥Python ◆開発 ○PHP" },
    @{ B = "【医療機関向け】既存システム改修・機能拡充エンジニア(Node.js/TypeScript)"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408214"; G = 185; H = "🔥TypeScript ◆Node.js" },
    @{ B = "【急募】携帯アプリ開発のプロフェッショナルを探しています!"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408519"; G = 100; H = "◆開発 ◇アプリ" },
    @{ B = "システムの開発補助や運営サポート【フルリモート×長期】"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408664"; G = 90;  H = "◆開発" },
    @{ B = "Javaプログラミング研修の演習サポート講師業務【経験不問】"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408635"; G = 85;  H = $null },
    @{ B = "2026年度新入社員研修Javaサブ講師"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408524"; G = 85;  H = "★Java" },
    @{ B = "2026年度新入社員研修Javaサブ講師 (4~6月)"; C = "システム開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408522"; G = 85;  H = "★Java" },
    @{ B = "2026年度新入社員研修Javaメイン講師"; C = "システム開発"; D = "1,000,000 円 ~ 3,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408503"; G = 85;  H = "★Java" },
    @{ B = "【急募】LLMによるMCP(Model Context Protocol)でのExcel操作機能開発"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5400689"; G = 75;  H = "◆開発" },
    @{ B = "仮想通貨トレードの運用とコンサル【1名】のみ募集"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5393695"; G = 55;  H = "◆コンサル" },
    @{ B = "【急募】社内Webアプリケーションのセキュリティ診断依頼"; C = "システム開発"; D = "500,000 円 ~ 1,000,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408563"; G = 45;  H = "◇アプリ" },
    @{ B = "クラウド(AWS/Azure) 運用管理 研修の演習サポート講師業務【経験不問】"; C = "システム開発"; D = "200,000 円 ~ 300,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408637"; G = 38;  H = "◇管理" },
    @{ B = "Notion実装・運用管理パートナー募集(長期・リモート)"; C = "システム開発"; D = "50,000 円 ~ 100,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408286"; G = 38;  H = "◇管理" },
    @{ B = "初回 iOSとAndroidのアプリ 課金(サブスク)"; C = "システム開発"; D = "10,000 円 ~ 20,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5398382"; G = 30;  H = "◇アプリ" },
    @{ B = "Ctrlを押しながらフォルダの上をクリックすると別窓で上の階層のフォルダが開くアプリの作成"; C = "システム開発"; D = "10,000 円 ~ 20,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408148"; G = 30;  H = "◇アプリ" },
    @{ B = "【急募】Oracleを活用したQ&Aシステム構築の依頼"; C = "システム開発"; D = "300,000 円 ~ 500,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408554"; G = 40;  H = $null },
    @{ B = "【急募】Teams Roomsの設定設置と保守サポート依頼"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408814"; G = 18;  H = $null },
    @{ B = "【急募】16タイプ診断コンテンツのLP制作"; C = "システム開発"; D = "100,000 円 ~ 200,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408735"; G = 18;  H = $null },
    @{ B = "bubbleでのサービス構築(difyとの連結)"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408871"; G = 13;  H = $null },
    @{ B = "急募 限定公開 限定公開の仕事"; C = "システム開発"; D = "20,000 円 ~ 50,000 円 / 固定"; E = "期限情報なし"; F = "https://www.lancers.jp/work/detail/5408546"; G = 13;  H = $null }
)

# Drop every existing hyperlink relationship up front; we rebuild them below so
# the F-column link targets always line up with the row they end up on.
$ws.Hyperlinks.Delete()

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $timestamp
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G

    if ($row.H -ne $null) {
        $ws.Cells.Item($r, 8).Value = $row.H
    } else {
        $ws.Cells.Item($r, 8).Value = ""
    }

    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $row.F)

    $r = $r + 1
}
